# Generate Report for handoff
# Renames the handed-off source doc id from
#   386960c8-7cc8-40c4-88a0-fa721e7d389c
# to
#   2d66d2e8-376a-4ce2-81f8-b74ce41da6fc
# and refreshes the handoff file hashes / timestamps that go with the new
# handoff round, across the Overview / zh-cn / de-de sheets. Hyperlink
# target URLs (the git blob history) are left untouched - only the cell
# text and the hyperlinks' visible display text change.

$wb = $excel.ActiveWorkbook

$oldId = "386960c8-7cc8-40c4-88a0-fa721e7d389c"
$newId = "2d66d2e8-376a-4ce2-81f8-b74ce41da6fc"

$oldZhHash = "e096ab349600b832ccdde027d8c187a3b4398cc8"
$newZhHash = "efd353867975518c7a40b5b4942a1754d71c0396"

$oldDeHash = $oldZhHash
$newDeHash = $newZhHash

$newMdName = "$newId.md"
$newZhXlf  = "$newId.$newZhHash.zh-cn.xlf"
$newDeXlf  = "$newId.$newDeHash.de-de.xlf"

$newZhDatetime = "2016-01-07 08:13:07"
$newDeDatetime = "2016-01-07 08:13:18"

# Hyperlink TARGET URLs are untouched by this edit (only the cell text /
# displayed link text changes) - keep them byte-identical to the
# originals, which still embed the *old* id / hash from the git history.
$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/5c82fd74d56cb8789224391be302a1817663aca6/e2e/$oldId.md"
$configAddress = "https://github.com/OpenLocalizationTest/oltest/blob/c8f1b1579cf85e747fdd5029ed41b259f4b72bdc/.localization-config"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/36b1e1255ce9ce0abc9521bb0afac32d3fcfcabd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$oldId.$oldZhHash.zh-cn.xlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3ffb0797177ce7efb65887fb54f87cde7919132/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$oldId.$oldDeHash.de-de.xlf"

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMdName) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configAddress, "", "", ".localization-config") | Out-Null

# ---- Sheet "zh-cn" ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("C2").Value = $newZhXlf
$wsZh.Range("D2").Value = $newZhDatetime

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddress, "", "", $newMdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfAddress, "", "", $newZhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configAddress, "", "", ".localization-config") | Out-Null

# ---- Sheet "de-de" ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("C2").Value = $newDeXlf
$wsDe.Range("D2").Value = $newDeDatetime

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddress, "", "", $newMdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfAddress, "", "", $newDeXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configAddress, "", "", ".localization-config") | Out-Null
